$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'61.305.29"
$ws.Range("E2").Value = "  +0.66%  "
$ws.Range("D3").Value = "'2.929.04"
$ws.Range("E3").Value = "  +0.45%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'594.04"
$ws.Range("E5").Value = "  +0.25%  "
$ws.Range("D6").Value = "'143.32"
$ws.Range("E6").Value = "  -1.28%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  -1.18%  "
$ws.Range("E9").Value = "  +1.39%  "
$ws.Range("E10").Value = "  -1.66%  "
$ws.Range("E11").Value = "  -0.40%  "
$ws.Range("E12").Value = "  -0.80%  "
$ws.Range("E13").Value = "  -0.97%  "
$ws.Range("E14").Value = "  +0.16%  "
$ws.Range("D15").Value = "'3.414.94"
$ws.Range("E15").Value = "  +0.49%  "
$ws.Range("D16").Value = "'61.294.18"
$ws.Range("E16").Value = "  +0.64%  "
$ws.Range("D17").Value = "'2.929.62"
$ws.Range("E17").Value = "  +0.45%  "
$ws.Range("D18").Value = "'6.63"
$ws.Range("E18").Value = "  -0.73%  "
$ws.Range("D19").Value = "'433.48"
$ws.Range("E19").Value = "  +0.90%  "
$ws.Range("D20").Value = "'13.53"
$ws.Range("E20").Value = "  +1.27%  "
$ws.Range("E21").Value = "  -1.02%  "
$ws.Range("D22").Value = "'7.05"
$ws.Range("E22").Value = "  -0.15%  "
$ws.Range("D23").Value = "'81.49"
$ws.Range("E23").Value = "  +0.12%  "
$ws.Range("D24").Value = "'10.80"
$ws.Range("E24").Value = "  -0.94%  "
$ws.Range("E25").Value = "  -1.98%  "
$ws.Range("D26").Value = "'11.70"
$ws.Range("E26").Value = "  -1.99%  "
$ws.Range("E27").Value = "  +0.05%  "
$ws.Range("E28").Value = "  -3.60%  "
$ws.Range("D29").Value = "'2.59"
$ws.Range("E29").Value = "  -0.98%  "
$ws.Range("D30").Value = "'6.87"
$ws.Range("E30").Value = "  -2.31%  "
$ws.Range("E31").Value = "  +0.95%  "
$ws.Range("D32").Value = "'0.108"
$ws.Range("E32").Value = "  +1.46%  "
$ws.Range("E33").Value = "  +0.00%  "
$ws.Range("D34").Value = "'0.0₃0871"
$ws.Range("E34").Value = "  +2.92%  "
$ws.Range("E35").Value = "  -0.30%  "
$ws.Range("D36").Value = "'5.62"
$ws.Range("E36").Value = "  -0.08%  "
$ws.Range("E37").Value = "  -1.99%  "
$ws.Range("E38").Value = "  +0.34%  "
$ws.Range("E39").Value = "  +0.35%  "
$ws.Range("D40").Value = "'8.49"
$ws.Range("E40").Value = "  -0.22%  "
$ws.Range("D41").Value = "'41.97"
$ws.Range("E41").Value = "  +5.27%  "
$ws.Range("E42").Value = "  -3.42%  "
$ws.Range("D43").Value = "'0.0343"
$ws.Range("E43").Value = "  -0.19%  "
$ws.Range("D44").Value = "'2.695.84"
$ws.Range("E44").Value = "  -0.10%  "
$ws.Range("D45").Value = "'133.56"
$ws.Range("E45").Value = "  +1.77%  "
$ws.Range("D46").Value = "'361.99"
$ws.Range("E46").Value = "  -3.47%  "
$ws.Range("E47").Value = "  +0.05%  "
$ws.Range("D48").Value = "'23.49"
$ws.Range("E48").Value = "  -1.55%  "
$ws.Range("E49").Value = "  -1.21%  "
$ws.Range("D50").Value = "'2.00"
$ws.Range("E50").Value = "  -0.69%  "
$ws.Range("D51").Value = "'0.125"
$ws.Range("E51").Value = "  -0.16%  "
